# Update gh-pages to output generated at 456a3b4
# Applies numeric "want-to-go count" bumps and a name-text fix (鸳 -> 鸢)
# across the "展览", "演出" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3109
$ws.Range("F3").Value = 514
$ws.Range("F5").Value = 73
$ws.Range("F6").Value = 18
$ws.Range("F8").Value = 24
$ws.Range("F10").Value = 15284
$ws.Range("F12").Value = 156
$ws.Range("F14").Value = 6043
$ws.Range("F18").Value = 99
$ws.Range("F21").Value = 109
$ws.Range("C25").Value = "苏州·代号鸢only茶话会-星渡咖啡"
$ws.Range("F26").Value = 4972
$ws.Range("F27").Value = 132
$ws.Range("F28").Value = 10905

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 18

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3109
$ws.Range("F4").Value = 514
$ws.Range("F6").Value = 73
$ws.Range("F7").Value = 18
$ws.Range("F9").Value = 24
$ws.Range("F11").Value = 15284
$ws.Range("F13").Value = 156
$ws.Range("F15").Value = 6043
$ws.Range("F19").Value = 99
$ws.Range("F22").Value = 109
$ws.Range("C26").Value = "苏州·代号鸢only茶话会-星渡咖啡"
$ws.Range("F27").Value = 4972
$ws.Range("F28").Value = 132
$ws.Range("F29").Value = 18
$ws.Range("F30").Value = 10905
